$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Wat ga ik doen?" (column D) entries for each team member.
# Setting these in this particular order reproduces the same shared-string
# table ordering that Excel itself produced when it rewrote the file.
$ws.Range("D6").Value = "scrum+to do lijst - trigger + stored procedures, festivals aanpassen,groepen aanpassen"
$ws.Range("D2").Value = "template maken, details groepen"
$ws.Range("D3").Value = "remote van git,home"
$ws.Range("D4").Value = "template maken, festival detail +festival detail aanpassen,groepen detail aanpassen"
$ws.Range("D5").Value = "template maken,festivals, groepen"

# Give column E an explicit width (matches the author's added <col> entry
# for the as-yet-unused next column).
$ws.Range("E1").ColumnWidth = 21.29

# Move / leave the active selection on D6, the last cell edited.
$ws.Range("D6").Select()
